# Scenario 6.xlsx edit
# - Decrease the "probability of fetal death before 4 weeks from conception"
#   values on the Phase1 and Phase2 sheets (rows 2-5, columns B/C respectively).
#   The dependent "remaining probability" formulas in D (Phase1) / E (Phase2)
#   recalculate automatically.
# - Update the saved selection / active-sheet state: Phase1 becomes the
#   active tab with B2:B5 selected, and Phase2's selection moves to C2:C5
#   (no longer the active tab).

$wb = $excel.ActiveWorkbook

$wsPhase1 = $wb.Worksheets.Item("Phase1")
$wsPhase2 = $wb.Worksheets.Item("Phase2")

# --- Update probability values -------------------------------------------------

# Phase1: column B holds "Prob Fetal Death", column D = 1 - C - B (formula already present)
$wsPhase1.Range("B2").Value = 0.1
$wsPhase1.Range("B3").Value = 0.1
$wsPhase1.Range("B4").Value = 0.05
$wsPhase1.Range("B5").Value = 0.05

# Phase2: column C holds "Prob Fetal Death", column E = 1 - D - C (formula already present)
$wsPhase2.Range("C2").Value = 0.1
$wsPhase2.Range("C3").Value = 0.1
$wsPhase2.Range("C4").Value = 0.05
$wsPhase2.Range("C5").Value = 0.05

# --- Update selections / active sheet ------------------------------------------

# Move Phase2's saved selection to C2:C5 (it is no longer the active tab).
[void]$wsPhase2.Range("C2:C5").Select()

# Make Phase1 the active sheet and select B2:B5.
[void]$wsPhase1.Activate()
[void]$wsPhase1.Range("B2:B5").Select()

$wb.Application.Calculate()
